$d = $word.ActiveDocument

# 1. "Our method can be applied" -> "Our method is an overdetermined inverse problem that can be applied"
$d.Content.Find.Execute(
    "Our method can be applied",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Our method is an overdetermined inverse problem that can be applied",
    2) | Out-Null

# 2. "with different magnetization" -> "with different but homogeneous magnetization"
$d.Content.Find.Execute(
    "with different magnetization",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "with different but homogeneous magnetization",
    2) | Out-Null

# 3. "direction of that as the ones in" -> "direction of that ones in"
$d.Content.Find.Execute(
    "direction of that as the ones in",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "direction of that ones in",
    2) | Out-Null

# 4. "strongly  suggesting" (double space) -> "strongly suggesting" (single space)
$d.Content.Find.Execute(
    "strongly  suggesting",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "strongly suggesting",
    2) | Out-Null
